$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper cells: some D-column values look like plain numbers
# (e.g. "0.999", "226.91") and Excel would auto-convert them to the Number
# type on assignment. The source data keeps these as text, matching the
# original inline-string cells, so we briefly mark the cell as Text format,
# assign the literal string, then restore the default "Normal" style so no
# stray number-format style is left behind on the cell.

$ws.Range("D2").Value = "37.190.21"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.027.60"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  -5.66%  "
$ws.Range("D12").Value = "2.316.45"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.744"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "2.039.40"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "37.145.38"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +6.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0216"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.471.76"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0914"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "2.209.35"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.46%  "
